# Update "想去人数" (number of people interested) values for a few events
# in both the "展览" and "全部类型" sheets, reflecting refreshed counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6041
$ws1.Range("F10").Value = 61
$ws1.Range("F14").Value = 438
$ws1.Range("F18").Value = 1673

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6041
$ws4.Range("F11").Value = 61
$ws4.Range("F15").Value = 438
$ws4.Range("F19").Value = 1673
